# Updates nombre_aides (C), nombre_entreprises (D) and montant_total (E)
# figures for the "Fonds de solidarite - volet 1" regional/categorie-juridique
# dataset, reflecting the 2022-05-12 data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "C10" = 345532
    "E10" = 1817701330
    "C11" = 1265
    "E11" = 46515184
    "C13" = 187832
    "E13" = 1165163461
    "C88" = 71261
    "E88" = 110287231
    "C91" = 18845
    "E91" = 75099352
    "C93" = 16917
    "E93" = 50434566
    "C98" = 6285
    "E98" = 19290013
    "C100" = 9333
    "E100" = 23704653
    "C112" = 145226
    "E112" = 716158923
    "C119" = 8982
    "E119" = 37087784
    "C121" = 1306111
    "E121" = 2274508151
    "C129" = 633318
    "E129" = 3426201653
    "C130" = 4239
    "E130" = 140350460
    "C132" = 585588
    "D132" = 90777
    "E132" = 3459563006
    "C139" = 76637
    "E139" = 114129393
    "C144" = 25065
    "E144" = 92331999
    "C145" = 72
    "E145" = 6534296
    "C146" = 7438
    "E146" = 37681278
    "C150" = 895
    "E150" = 2015502
    "C151" = 39920
    "E151" = 60358891
    "C154" = 18430
    "E154" = 72554723
    "C156" = 12394
    "E156" = 40027556
    "C186" = 236818
    "E186" = 1189702070
    "C194" = 18378
    "E194" = 71345179
    "C215" = 230252
    "E215" = 408700746
    "C221" = 135494
    "E221" = 681816845
    "C240" = 205897
    "E240" = 1068613038
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

